$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old data area (rows 2-20, columns A-I)
$ws.Range("A2:I20").ClearContents()

# Scratch cell used to force text (non-numeric) typing for numeric-looking strings
# (dates in column A, id_venda numbers in column D) without leaving stray styles
# on the destination cells themselves.
$scratch = $ws.Range("K1")
$scratch.NumberFormat = "@"

function Set-TextValue($cell, $text) {
    $scratch.Value = $text
    $scratch.Copy()
    $cell.PasteSpecial(-4163)
}

$data = @(
    @('2025-07-14', 250, 'ERAM ESTALEIRO RIO AMAZONAS LTDA', '54263882', 15011531, 'DETERGENTE LIMPOL COCO 500ML', 133, 21.89, 34.09),
    @('2025-07-14', 60, 'SAWEM DA AMAZONIA LTDA', '54285324', 27262762, 'MARCA TEXTO AMARELO UND JOCAR OFFICE', 104, 7.32, 8.15),
    @('2025-07-15', 130, 'MUSASHI DA AMAZONIA LTDA', '54346779', 17171383, 'DETERGENTE DESENGRAX MAX PINE AUDAX 5L', -9, 13.88, 29.36),
    @('2025-07-15', 300, 'MUSASHI DA AMAZONIA LTDA', '54346779', 14795919, 'SACO DE LIXO 200L COMUM PACOTINHO C/5 UND SACOLMAX', 244, 43.84, 61.77),
    @('2025-07-16', 70, 'TECHLOG - SERVICOS DE GESTAO E SISTEMAS INFORMATIZ', '54316914', 19264853, 'PAPEL TOALHA INTERFOLHADO 1250 FLS NEWPAPER 100% CELULOSE', 722, 8.22, 13.46),
    @('2025-07-16', 10, 'V V REFEICOES LTDA', '54396269', 17125814, 'COADOR DE CAFE P', 21, 2.17, 2.01),
    @('2025-07-16', 40, 'CONDOMINIO RESIDENCIAL EPHYGENIO SALLES', '54365832', 20619556, 'PEDRA SANITARIA 35G RUBI FLORAL', 136, 13.17, 8.54),
    @('2025-07-16', 10, 'V V REFEICOES LTDA', '54396269', 17541022, 'COADOR DE CAFE INDUSTRIAL (MAIOR)', 26, 2.09, 2.05),
    @('2025-07-16', 10, 'CONDOMINIO DO TVLANDIA MALL', '54119372', 28133466, 'DESINFETANTE CONCENTRADO 5L AUDAX MAX 1:200 - LAVANDA', 105, 2.37, 1.91),
    @('2025-07-16', 60, 'V V REFEICOES LTDA', '54396269', 42173656, 'COADOR DE CAFÉ G', 23, 6.07, 15.64),
    @('2025-07-17', 100, 'TEL TELECOMUNICACOES LTDA.', '54443314', 11936640, 'LIMPADOR VEJA MULTIUSO GOLD 500ML', 3364, 10.53, 14.59),
    @('2025-07-21', 240, 'JURUA ESTALEIROS E NAVEGACAO LTDA', '54563773', 16933123, 'SABAO EM PO LAVAGEM PERFEITA OMO  - 400G', 67, 33.69, 62.9),
    @('2025-07-21', 200, 'JURUA ESTALEIROS E NAVEGACAO LTDA', '54503121', 13996941, 'SACO DE LIXO 30L REFORCADO PACOTINHO C/10 UND FORTE MAX', 81, 13.55, 22.47),
    @('2025-07-22', 22, 'CONDOMINIO CRISTAL TOWER', '54627290', 11939645, 'SACO DE LIXO 50L PRETO COMUM - PCT C/100 UND', 154, 4.3, 4.78),
    @('2025-07-22', 9, 'CARITAS ARQUIDIOCESANA DE MANAUS', '54625975', 17059594, 'SABAO EM PO ESPUMIL 4KG', 1, 1.94, 1.87),
    @('2025-07-23', 400, 'V V REFEICOES LTDA', '54692772', 32130390, 'ESPONJA MULTIUSO JEITOSA', 11689, 21.45, 51.21),
    @('2025-07-24', 40, 'MUNDY SERVIS LIMPEZAS PREDIAIS LTDA', '54704921', 20619556, 'PEDRA SANITARIA 35G RUBI FLORAL', 136, 13.17, 8.54),
    @('2025-07-24', 70, 'SUPER G - TRANSPORTE E LOCACAO DE EQUIPAMENTOS INDUSTRIAIS L', '54736983', 12097805, 'FLANELA LARANJA TAM G 38x58CM', 24, 15.16, 13.57),
    @('2025-07-25', 160, 'JURUA ESTALEIROS E NAVEGACAO LTDA', '54791493', 11939672, 'VASSOURA PIACAVA 20 FUROS', 298, 10.05, 29.14),
    @('2025-07-25', 100, 'MANJAR SERVICOS GERAIS SA', '54799813', 31186309, 'LA DE ACO ASSOLAN 45G', 28, 9.59, 17.49),
    @('2025-07-25', 168, 'JURUA ESTALEIROS E NAVEGACAO LTDA', '54791493', 15011531, 'DETERGENTE LIMPOL COCO 500ML', 133, 21.89, 34.09)
)

$r = 2
foreach ($row in $data) {
    Set-TextValue $ws.Cells.Item($r, 1) $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    Set-TextValue $ws.Cells.Item($r, 4) $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $r = $r + 1
}

# Clean up scratch column so it does not linger in the saved worksheet
$scratch.Clear()
$ws.Columns("K").Delete()

Write-Host "Rows written:" ($r - 2)
